$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 18.21482920925333
$ws.Range("C2").Value = 16.87178653249114
$ws.Range("D2").Value = 15.38138818166248
$ws.Range("E2").Value = 16.84162632776718
$ws.Range("G2").Value = 3.782032829084129
$ws.Range("I2").Value = 36.65488141297421
$ws.Range("J2").Value = 9.654571020340018
$ws.Range("K2").Value = 19.97528673912249
$ws.Range("N2").Value = 25.21721950639624
# Row 3
$ws.Range("B3").Value = 18.05593116864422
$ws.Range("C3").Value = 16.72386369035351
$ws.Range("D3").Value = 15.33660753858612
$ws.Range("E3").Value = 16.79827313350749
$ws.Range("G3").Value = 3.785915127842101
$ws.Range("I3").Value = 36.58803499910546
$ws.Range("J3").Value = 9.668366851138012
$ws.Range("K3").Value = 19.86367723556559
$ws.Range("N3").Value = 25.22639465723728
# Row 4
$ws.Range("B4").Value = 17.9631415513132
$ws.Range("C4").Value = 16.63710699999696
$ws.Range("D4").Value = 15.31270764385304
$ws.Range("E4").Value = 16.77560865704005
$ws.Range("G4").Value = 3.788420884580205
$ws.Range("I4").Value = 36.55258076001053
$ws.Range("J4").Value = 9.678372066819971
$ws.Range("K4").Value = 19.80020782177519
$ws.Range("N4").Value = 25.23377773216287
# Row 5
$ws.Range("B5").Value = 17.92656747170013
$ws.Range("C5").Value = 16.60280802469884
$ws.Range("D5").Value = 15.3038774589037
$ws.Range("E5").Value = 16.76737152795427
$ws.Range("G5").Value = 3.789472802370455
$ws.Range("I5").Value = 36.53954215439561
$ws.Range("J5").Value = 9.682834785597349
$ws.Range("K5").Value = 19.77563543541691
$ws.Range("N5").Value = 25.23722579262691
# Row 6
$ws.Range("B6").Value = 17.92057022819341
$ws.Range("C6").Value = 16.59717732183994
$ws.Range("D6").Value = 15.30246625672595
$ws.Range("E6").Value = 16.76606419497199
$ws.Range("G6").Value = 3.789649336416142
$ws.Range("I6").Value = 36.53746231011901
$ws.Range("J6").Value = 9.683599089196223
$ws.Range("K6").Value = 19.77163377286753
$ws.Range("N6").Value = 25.23782486277911
# Row 7
$ws.Range("B7").Value = 17.96264323804917
$ws.Range("C7").Value = 16.6366401194423
$ws.Range("D7").Value = 15.31258486954207
$ws.Range("E7").Value = 16.77549351892898
$ws.Range("G7").Value = 3.788434946238032
$ws.Range("I7").Value = 36.55239920671188
$ws.Range("J7").Value = 9.678430692199157
$ws.Range("K7").Value = 19.79987117465119
$ws.Range("N7").Value = 25.23382245564324
# Row 8
$ws.Range("B8").Value = 18.15907373384744
$ws.Range("C8").Value = 16.81995818474225
$ws.Range("D8").Value = 15.36520488340018
$ws.Range("E8").Value = 16.82586033128991
$ws.Range("G8").Value = 3.783346195769147
$ws.Range("I8").Value = 36.63067511343252
$ws.Range("J8").Value = 9.65900913361428
$ws.Range("K8").Value = 19.93576618893931
$ws.Range("N8").Value = 25.22001967951944
# Row 9
$ws.Range("B9").Value = 18.58037328062182
$ws.Range("C9").Value = 17.21021877049724
$ws.Range("D9").Value = 15.49665650008066
$ws.Range("E9").Value = 16.95576907974403
$ws.Range("G9").Value = 3.774329728077898
$ws.Range("I9").Value = 36.82832875988306
$ws.Range("J9").Value = 9.633113982767945
$ws.Range("K9").Value = 20.24143693769108
$ws.Range("N9").Value = 25.20685845182338
# Row 10
$ws.Range("B10").Value = 18.90937579168101
$ws.Range("C10").Value = 17.51356072722907
$ws.Range("D10").Value = 15.61004829658843
$ws.Range("E10").Value = 17.06981528096042
$ws.Range("G10").Value = 3.768284377862504
$ws.Range("I10").Value = 37.00014894870006
$ws.Range("J10").Value = 9.621538890116522
$ws.Range("K10").Value = 20.48848220489884
$ws.Range("N10").Value = 25.20570131014714
# Row 11
$ws.Range("B11").Value = 19.06270335609913
$ws.Range("C11").Value = 17.65468174854854
$ws.Range("D11").Value = 15.66517619751905
$ws.Range("E11").Value = 17.12563554355376
$ws.Range("G11").Value = 3.765658268119881
$ws.Range("I11").Value = 37.08400469111787
$ws.Range("J11").Value = 9.617893627224552
$ws.Range("K11").Value = 20.60541114537334
$ws.Range("N11").Value = 25.20703023283323
# Row 12
$ws.Range("B12").Value = 19.12124156851994
$ws.Range("C12").Value = 17.70852838977692
$ws.Range("D12").Value = 15.68655095581659
$ws.Range("E12").Value = 17.14732966466502
$ws.Range("G12").Value = 3.764681524874727
$ws.Range("I12").Value = 37.11656830985655
$ws.Range("J12").Value = 9.616746358656043
$ws.Range("K12").Value = 20.65031230347618
$ws.Range("N12").Value = 25.20780069050338
# Row 13
$ws.Range("B13").Value = 19.10861393637639
$ws.Range("C13").Value = 17.69691413037727
$ws.Range("D13").Value = 15.6819255060585
$ws.Range("E13").Value = 17.14263288516664
$ws.Range("G13").Value = 3.764891098311321
$ws.Range("I13").Value = 37.10951932075757
$ws.Range("J13").Value = 9.616983074563231
$ws.Range("K13").Value = 20.64061484894227
$ws.Range("N13").Value = 25.20762286736669
# Row 14
$ws.Range("B14").Value = 19.06751007264634
$ws.Range("C14").Value = 17.65910383093165
$ws.Range("D14").Value = 15.66692477551486
$ws.Range("E14").Value = 17.12740924751523
$ws.Range("G14").Value = 3.765577556645996
$ws.Range("I14").Value = 37.0866675584701
$ws.Range("J14").Value = 9.617794568386262
$ws.Range("K14").Value = 20.60909289326363
$ws.Range("N14").Value = 25.20708826097592
# Row 15
$ws.Range("B15").Value = 19.04239330650104
$ws.Range("C15").Value = 17.63599567661734
$ws.Range("D15").Value = 15.65780104696983
$ws.Range("E15").Value = 17.11815644012267
$ws.Range("G15").Value = 3.766000334790057
$ws.Range("I15").Value = 37.07277528675528
$ws.Range("J15").Value = 9.61832199231222
$ws.Range("K15").Value = 20.58986493557549
$ws.Range("N15").Value = 25.20679561158915
# Row 16
$ws.Range("B16").Value = 18.89942489246456
$ws.Range("C16").Value = 17.50439745946949
$ws.Range("D16").Value = 15.60651599664892
$ws.Range("E16").Value = 17.06624571329687
$ws.Range("G16").Value = 3.768458484550566
$ws.Range("I16").Value = 36.99478256334266
$ws.Range("J16").Value = 9.621809732694862
$ws.Range("K16").Value = 20.48092948531387
$ws.Range("N16").Value = 25.20565183342427
# Row 17
$ws.Range("B17").Value = 18.81262093332167
$ws.Range("C17").Value = 17.42443743579524
$ws.Range("D17").Value = 15.57595492247249
$ws.Range("E17").Value = 17.03540221517571
$ws.Range("G17").Value = 3.769998143544237
$ws.Range("I17").Value = 36.9483884324908
$ws.Range("J17").Value = 9.624364450720483
$ws.Range("K17").Value = 20.41524392289504
$ws.Range("N17").Value = 25.20542567311637
# Row 18
$ws.Range("B18").Value = 18.76304140374364
$ws.Range("C18").Value = 17.37874369517817
$ws.Range("D18").Value = 15.55871141300672
$ws.Range("E18").Value = 17.01803338939821
$ws.Range("G18").Value = 3.770895388411295
$ws.Range("I18").Value = 36.92224032019079
$ws.Range("J18").Value = 9.62598636505883
$ws.Range("K18").Value = 20.37789324023684
$ws.Range("N18").Value = 25.20547019938722
# Row 19
$ws.Range("B19").Value = 18.74631588536912
$ws.Range("C19").Value = 17.36332490962587
$ws.Range("D19").Value = 15.55293080800481
$ws.Range("E19").Value = 17.01221671690245
$ws.Range("G19").Value = 3.771201188632593
$ws.Range("I19").Value = 36.91347944804038
$ws.Range("J19").Value = 9.626561706170481
$ws.Range("K19").Value = 20.36532170680681
$ws.Range("N19").Value = 25.20551525035265
# Row 20
$ws.Range("B20").Value = 18.82182572808228
$ws.Range("C20").Value = 17.43291887112332
$ws.Range("D20").Value = 15.57917366623128
$ws.Range("E20").Value = 17.03864717770083
$ws.Range("G20").Value = 3.769833036872953
$ws.Range("I20").Value = 36.95327170036322
$ws.Range("J20").Value = 9.624076711571391
$ws.Range("K20").Value = 20.42219199002972
$ws.Range("N20").Value = 25.20543167377951
# Row 21
$ws.Range("B21").Value = 19.0795707484038
$ws.Range("C21").Value = 17.67019892574371
$ws.Range("D21").Value = 15.67131740070928
$ws.Range("E21").Value = 17.1318657909773
$ws.Range("G21").Value = 3.765375447541582
$ws.Range("I21").Value = 37.09335779206308
$ws.Range("J21").Value = 9.617549885760869
$ws.Range("K21").Value = 20.61833500968096
$ws.Range("N21").Value = 25.20723803222559
# Row 22
$ws.Range("B22").Value = 19.25077512385904
$ws.Range("C22").Value = 17.82762829194676
$ws.Range("D22").Value = 15.73444173400511
$ws.Range("E22").Value = 17.19602606400703
$ws.Range("G22").Value = 3.762565314435658
$ws.Range("I22").Value = 37.18962412268632
$ws.Range("J22").Value = 9.614642997016265
$ws.Range("K22").Value = 20.75014017326327
$ws.Range("N22").Value = 25.20997634522753
# Row 23
$ws.Range("B23").Value = 19.15916496808055
$ws.Range("C23").Value = 17.74340423527483
$ws.Range("D23").Value = 15.70048923931082
$ws.Range("E23").Value = 17.16149003395771
$ws.Range("G23").Value = 3.764055734651881
$ws.Range("I23").Value = 37.13781720618915
$ws.Range("J23").Value = 9.61607010994395
$ws.Range("K23").Value = 20.67947314605815
$ws.Range("N23").Value = 25.20837218971679
# Row 24
$ws.Range("B24").Value = 18.81766322356749
$ws.Range("C24").Value = 17.42908355022783
$ws.Range("D24").Value = 15.57771745452068
$ws.Range("E24").Value = 17.03717899710762
$ws.Range("G24").Value = 3.769907643981042
$ws.Range("I24").Value = 36.95106234057485
$ws.Range("J24").Value = 9.624206321350803
$ws.Range("K24").Value = 20.41904948068415
$ws.Range("N24").Value = 25.2054284172323
# Row 25
$ws.Range("B25").Value = 18.46279040953723
$ws.Range("C25").Value = 17.10156348839392
$ws.Range("D25").Value = 15.45810868568681
$ws.Range("E25").Value = 16.91732700269079
$ws.Range("G25").Value = 3.77666667593426
$ws.Range("I25").Value = 36.77015960462588
$ws.Range("J25").Value = 9.638812089434253
$ws.Range("K25").Value = 20.15468950237308
$ws.Range("N25").Value = 25.20892641474257
